# Big stimulus update:
#  - rename the "face" stimulus category to "book" (face//face_NN.jpg -> book//book_NN.jpg)
#  - expand the abbreviated correct_ans codes (y/r/b) to their full words (left/right/center)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange

$firstRow = $ur.Row
$firstCol = $ur.Column
$rows = $ur.Rows.Count
$cols = $ur.Columns.Count

# Pull the whole used range into memory in one shot (1-based array, [row, col])
$vals = $ur.Value()

# Locate the "correct_ans" column from the header row instead of hard-coding it,
# so the script is robust to column order.
$correctAnsCol = -1
for ($c = 1; $c -le $cols; $c++) {
  $header = $vals[1, $c]
  if ($header -eq "correct_ans") {
    $correctAnsCol = $c
  }
}

# Map of abbreviation -> full word used in the correct_ans column
$ansMap = @{ "y" = "left"; "r" = "right"; "b" = "center" }

for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $v = $vals[$r, $c]

    if ($v -is [string]) {
      # Rename the "face" stimulus folder/filename prefix to "book" wherever it occurs
      if ($v.Contains("face//face_")) {
        $vals[$r, $c] = $v.Replace("face", "book")
        continue
      }

      # Expand the correct_ans abbreviation codes (skip the header row)
      if ($correctAnsCol -gt 0 -and $c -eq $correctAnsCol -and $r -gt 1) {
        if ($ansMap.ContainsKey($v)) {
          $vals[$r, $c] = $ansMap[$v]
        }
      }
    }
  }
}

$ur.Value = $vals
